$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 51: A Bile Business (item id 5486)
$ws.Range("H51").Value = 8869.143
$ws.Range("J51").Value = 10216.8
$ws.Range("L51").Value = 10216.8
$ws.Range("N51").Value = -11184.8
# Row 107: Another Man's Ink (item id 27766)
$ws.Range("H107").Value = 945.125
$ws.Range("J107").Value = 268.5
$ws.Range("L107").Value = 268.5
$ws.Range("N107").Value = -4108.5
# Row 127: Liquid Competence (item id 36114)
$ws.Range("H127").Value = 2215.8
$ws.Range("I127").Value = 1644.75
$ws.Range("J127").Value = 4500
$ws.Range("K127").Value = 4934.25
$ws.Range("L127").Value = 13500
$ws.Range("M127").Value = 25.75
$ws.Range("N127").Value = -23420
# Row 132: Fast-forwarding Flora (item id 44049)
$ws.Range("H132").Value = 4338.4443
$ws.Range("I132").Value = 3729.6667
$ws.Range("K132").Value = 11189.0001
$ws.Range("M132").Value = -8659.000100000001
# Row 137: Cutting Edge of Culinary Quality (item id 44013)
$ws.Range("H137").Value = 3736.3333
$ws.Range("I137").Value = 1449.1818
$ws.Range("J137").Value = 13799.8
$ws.Range("K137").Value = 4347.5454
$ws.Range("L137").Value = 41399.39999999999
$ws.Range("M137").Value = -1797.5454
$ws.Range("N137").Value = -46499.39999999999
# Row 138: All-night Crafting (item id 44169)
$ws.Range("H138").Value = 5542.473
$ws.Range("J138").Value = 4770.8413
$ws.Range("L138").Value = 14312.5239
$ws.Range("N138").Value = -24592.5239
# Row 139: Something Salty and Ceremonial (item id 42306)
$ws.Range("H139").Value = 97662
$ws.Range("J139").Value = 97662
$ws.Range("L139").Value = 97662
$ws.Range("N139").Value = -107942

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust (item id 44147)
$ws.Range("H32").Value = 19009.77
$ws.Range("I32").Value = 11853.116
$ws.Range("K32").Value = 11853.116
$ws.Range("M32").Value = -11566.116
# Row 74: As the Bolt Flies (item id 44000)
$ws.Range("H74").Value = 44617.5
$ws.Range("I74").Value = 50358.094
$ws.Range("J74").Value = 4433.3335
$ws.Range("K74").Value = 50358.094
$ws.Range("L74").Value = 4433.3335
$ws.Range("M74").Value = -49484.094
$ws.Range("N74").Value = -6181.3335
# Row 77: Heavy Metal Banned (L) (item id 44000)
$ws.Range("H77").Value = 44617.5
$ws.Range("I77").Value = 50358.094
$ws.Range("J77").Value = 4433.3335
$ws.Range("K77").Value = 251790.47
$ws.Range("L77").Value = 22166.6675
$ws.Range("M77").Value = -247422.47
$ws.Range("N77").Value = -30902.6675
# Row 132: Don't Bore Me, Ore Me (item id 43997)
$ws.Range("H132").Value = 80743.38
$ws.Range("I132").Value = 94696.73
$ws.Range("K132").Value = 284090.19
$ws.Range("M132").Value = -281560.19

$ws = $wb.Worksheets.Item("BSM")
# Row 21: Awl or Nothing (item id 19542)
$ws.Range("H21").Value = 18417.857
$ws.Range("J21").Value = 18417.857
$ws.Range("L21").Value = 18417.857
$ws.Range("N21").Value = -18889.857
# Row 94: High Steal (item id 19939)
$ws.Range("H94").Value = 8606.75
$ws.Range("I94").Value = 8847.200000000001
$ws.Range("K94").Value = 8847.200000000001
$ws.Range("M94").Value = -8396.200000000001
# Row 105: Ingot to Wing It (item id 19947)
$ws.Range("H105").Value = 1910.25
$ws.Range("I105").Value = 1599.5
$ws.Range("K105").Value = 1599.5
$ws.Range("M105").Value = 147.5
# Row 132: Always Be Prepaired (item id 41855)
$ws.Range("H132").Value = 100373.4
$ws.Range("J132").Value = 100373.4
$ws.Range("L132").Value = 100373.4
$ws.Range("N132").Value = -110493.4
# Row 134: Ruthenium Supremium (item id 43998)
$ws.Range("H134").Value = 2486.7778
$ws.Range("I134").Value = 2134.1333
$ws.Range("K134").Value = 6402.3999
$ws.Range("M134").Value = -3867.3999

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof (item id 27691)
$ws.Range("H16").Value = 1684.1428
$ws.Range("I16").Value = 1131.5
$ws.Range("K16").Value = 1131.5
$ws.Range("M16").Value = -844.5
# Row 22: Driving Up the Wall (item id 5367)
$ws.Range("H22").Value = 1179.9286
$ws.Range("I22").Value = 1508.625
$ws.Range("K22").Value = 1508.625
$ws.Range("M22").Value = -1158.625
# Row 105: Zelkova, My Love (item id 19928)
$ws.Range("H105").Value = 2181.4375
$ws.Range("I105").Value = 2146.6924
$ws.Range("J105").Value = 2332
$ws.Range("K105").Value = 2146.6924
$ws.Range("L105").Value = 2332
$ws.Range("M105").Value = -399.6923999999999
$ws.Range("N105").Value = -5826
# Row 113: Patient Patients (item id 27691)
$ws.Range("H113").Value = 1684.1428
$ws.Range("I113").Value = 1131.5
$ws.Range("K113").Value = 1131.5
$ws.Range("M113").Value = 1038.5
# Row 132: Hull Lotta Damage (item id 44019)
$ws.Range("H132").Value = 2518.2
$ws.Range("I132").Value = 2518.8276
$ws.Range("K132").Value = 7556.4828
$ws.Range("M132").Value = -5026.4828

$ws = $wb.Worksheets.Item("CUL")
# Row 7: It's Always Sunny in Vylbrand (item id 4728)
$ws.Range("H7").Value = 1663.25
$ws.Range("I7").Value = 399
$ws.Range("K7").Value = 1197
$ws.Range("M7").Value = -1085

$ws = $wb.Worksheets.Item("GSM")
# Row 59: Sew Not Doing This (item id 2453)
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
# Row 70: Sky Is the Limit (item id 14146)
$ws.Range("H70").Value = 4865.6665
$ws.Range("I70").Value = 4799
$ws.Range("K70").Value = 4799
$ws.Range("M70").Value = -4529
# Row 73: Hulls of Broken Dreams (L) (item id 14146)
$ws.Range("H73").Value = 4865.6665
$ws.Range("I73").Value = 4799
$ws.Range("K73").Value = 4799
$ws.Range("M73").Value = -3863
# Row 80: Needs More Prayerbell (item id 12521)
$ws.Range("H80").Value = 3166.1428
$ws.Range("J80").Value = 4164.3335
$ws.Range("L80").Value = 4164.3335
$ws.Range("N80").Value = -6160.3335
# Row 83: With a Noise That Reaches Heaven (L) (item id 12521)
$ws.Range("H83").Value = 3166.1428
$ws.Range("J83").Value = 4164.3335
$ws.Range("L83").Value = 20821.6675
$ws.Range("N83").Value = -30805.6675
# Row 132: On Board for Lar (item id 44008)
$ws.Range("H132").Value = 63418.332
$ws.Range("I132").Value = 80761
$ws.Range("J132").Value = 2719
$ws.Range("K132").Value = 242283
$ws.Range("L132").Value = 8157
$ws.Range("M132").Value = -239753
$ws.Range("N132").Value = -13217

$ws = $wb.Worksheets.Item("LTW")
# Row 40: Best Served Toad (item id 36248)
$ws.Range("H40").Value = 3799.7742
$ws.Range("I40").Value = 3617.6428
$ws.Range("K40").Value = 3617.6428
$ws.Range("M40").Value = -3481.6428
# Row 55: It's Not a Job, It's a Calling (item id 5284)
$ws.Range("H55").Value = 756.8570999999999
$ws.Range("I55").Value = 491.9375
$ws.Range("J55").Value = 1110.0834
$ws.Range("K55").Value = 491.9375
$ws.Range("L55").Value = 1110.0834
$ws.Range("M55").Value = -318.9375
$ws.Range("N55").Value = -1456.0834
# Row 132: Tenets of Tanning (item id 44058)
$ws.Range("H132").Value = 28950.191
$ws.Range("I132").Value = 36457.582
$ws.Range("J132").Value = 4380.5454
$ws.Range("K132").Value = 109372.746
$ws.Range("L132").Value = 13141.6362
$ws.Range("M132").Value = -106842.746
$ws.Range("N132").Value = -18201.6362

$ws = $wb.Worksheets.Item("WVR")
# Row 120: A Turban for the Ages (item id 26310)
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
# Row 126: A Polished Purchase (item id 36210)
$ws.Range("H126").Value = 54421.844
$ws.Range("I126").Value = 84985
$ws.Range("K126").Value = 254955
$ws.Range("M126").Value = -252485
# Row 132: Comfy Cabins (item id 44029)
$ws.Range("H132").Value = 20886.508
$ws.Range("I132").Value = 21063
$ws.Range("K132").Value = 63189
$ws.Range("M132").Value = -60659
# Row 136: Weaving the Envelope (item id 44031)
$ws.Range("H136").Value = 10557927
$ws.Range("I136").Value = 1224007
$ws.Range("K136").Value = 3672021
$ws.Range("M136").Value = -3669471
